$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (row 12) into the new row 13, A column
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)

# Fill in the new row of data (2021年)
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 15084
$ws.Range("C13").Value = 18407
$ws.Range("D13").Value = 12467.98
$ws.Range("E13").Value = 6743.27
$ws.Range("F13").Value = 14339
$ws.Range("G13").Value = 21.2086
$ws.Range("H13").Value = 23435.54
$ws.Range("I13").Value = 1548
$ws.Range("J13").Value = 16896.54
$ws.Range("K13").Value = 47936
$ws.Range("L13").Value = 1131.2841
$ws.Range("M13").Value = 21830
$ws.Range("N13").Value = 66642
$ws.Range("O13").Value = 37.8774
